# Generate Report for Handback
# Updates the handback-status report with the latest handoff/handback
# timestamps recorded for the de-de localization file
# (c53545e8-35e1-47da-980b-21a25fc16a1f), and refreshes the corresponding
# datetime cells on the zh-cn sheet as part of the same report generation.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row for df5fe72b... file, refresh handoff/handback datetimes ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-23 14:51:19"
$zhcn.Range("H3").Value = "2016-03-23 14:52:05"

# --- de-de sheet: row for c53545e8... file (row 2), new handback recorded ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "c53545e8-35e1-47da-980b-21a25fc16a1f.a420034b7170d0a49a6215d5e0555ee4d71f640a.de-de.xlf"
$dede.Range("E2").Value = "2016-03-23 14:53:18"
$dede.Range("G2").Value = "c53545e8-35e1-47da-980b-21a25fc16a1f.a420034b7170d0a49a6215d5e0555ee4d71f640a.de-de.xlf"
$dede.Range("H2").Value = "2016-03-23 14:53:50"

# --- de-de sheet: row for df5fe72b... file (row 3), refresh handoff/handback datetimes ---
$dede.Range("D3").Value = "df5fe72b-3f64-4d1a-95bb-2ba78bc8179b.8011e602df58a08fc86b047b188d6e3d53c7b052.de-de.xlf"
$dede.Range("E3").Value = "2016-03-23 14:51:25"
$dede.Range("G3").Value = "df5fe72b-3f64-4d1a-95bb-2ba78bc8179b.8011e602df58a08fc86b047b188d6e3d53c7b052.de-de.xlf"
$dede.Range("H3").Value = "2016-03-23 14:52:15"
